$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "Dr. Alshimaa Atef, Dr. Heba Mahmoud Ali, Dr. Hend Mahmoud, Dr. Servinaz Sayed Mohammad, Dr. Rana Abo-Zaid, Dr. Veronia Rafat, Dr. Shimaa Ahmad Mekki, Dr. Amira Sobhy"
$ws.Range("G3").Value = "Administrator, Dr. Gehan Adel, Dr. Alshimaa Atef, Dr. Manar Montaser"
$ws.Range("G4").Value = "Dr. Menna tuâ€™Allah Medhat, Dr. Asmaa Reda, Dr. Nourhan Mahmoud, Dr. Heba Mahmoud Ali, Dr. Majorelle Magdy, Dr. Hanan Ragab, Dr. Shimaa Ahmad Mekki"
$ws.Range("G9").Value = "Dr. Marina Youhanna, Dr. Yasmeena Fattoh, Dr. Madeha Saeed, Dr. Eman M. Abo-Sakaya"
$ws.Range("G10").Value = "Dr. Basma Hamed, Dr. Amira Ibrahim"
$ws.Range("G12").Value = "Dr. Mona Ibrahim Hussein, Dr. Heba Al-Sayed Mohammad, Dr. Dalia Tarek Elwan"
$ws.Range("G17").Value = "Dr. Walaa Ghanima, Dr. Marian Samir, Dr. Enas Omran"
$ws.Range("G18").Value = "Dr. Eman Samir Gabry, Dr. Wafaa Ebida, Dr. Ola Abd Al-Fattah, Dr. Abdullah El-Agrody"
$ws.Range("G19").Value = "Dr. Eman Samir Gabry, Dr. Yasmin, Dr. Wafaa Ebida, Dr. Neveen Nashaat, Dr. Marina Sorial"
$ws.Range("G20").Value = "Dr. Remon, Dr. Marina Atef, Dr. Nardine, Dr. Yasmin, Dr. Monica, Dr. Neveen Nashaat, Dr. Marina Sorial"
$ws.Range("G21").Value = "Dr. Alshimaa Atef, Dr. Heba Mahmoud Ali, Dr. Hend Mahmoud, Dr. Servinaz Sayed Mohammad, Dr. Rana Abo-Zaid, Dr. Veronia Rafat, Dr. Shimaa Ahmad Mekki, Dr. Amira Sobhy"
$ws.Range("G22").Value = "Administrator, Dr. Gehan Adel, Dr. Alshimaa Atef, Dr. Manar Montaser"
$ws.Range("G23").Value = "Dr. Menna tuâ€™Allah Medhat, Dr. Asmaa Reda, Dr. Nourhan Mahmoud, Dr. Heba Mahmoud Ali, Dr. Majorelle Magdy, Dr. Hanan Ragab, Dr. Shimaa Ahmad Mekki"
$ws.Range("G24").Value = "Dr. Amera Ahmad Saad, Dr. Lamiaa Ossama, Dr. Fatma Elhady, Dr. Nada Mohammad, Dr. Abeer Ragab"
$ws.Range("G28").Value = "Dr. Dina Adel, Dr. Basma Hamed, Dr. Nourhan Osama, Dr. Arwa Al-Sayed, Dr. Esraa Mostafa, Dr. Sarah Abdelmohsen, Dr. Madeha Saeed, Dr. Marwa Mustafa, Dr. Yasmeena Fattoh, Dr. Eman M. Abo-Sakaya"
$ws.Range("G29").Value = "Dr. Esraa Mostafa, Dr. Amira Ibrahim, Dr. Yasmeena Fattoh"
$ws.Range("G31").Value = "Dr. Mona Ibrahim Hussein, Dr. Heba Al-Sayed Mohammad, Dr. Dalia Tarek Elwan"
$ws.Range("G32").Value = "Menna tuâ€™Allah Gamil, Dr. Nouran Mahmoud"
$ws.Range("G36").Value = "Dr. Walaa Ghanima, Dr. Marian Samir, Dr. Enas Omran"
$ws.Range("G37").Value = "Dr. Eman Samir Gabry, Dr. Wafaa Ebida, Dr. Ola Abd Al-Fattah, Dr. Abdullah El-Agrody"
$ws.Range("G38").Value = "Dr. Remon, Dr. Marina Atef, Dr. Nardine, Dr. Yasmin, Dr. Monica, Dr. Neveen Nashaat, Dr. Marina Sorial"
$ws.Range("G39").Value = "Dr. Remon, Dr. Marina Atef, Dr. Nardine, Dr. Yasmin, Dr. Monica, Dr. Neveen Nashaat, Dr. Marina Sorial"
$ws.Range("G40").Value = "Dr. Alshimaa Atef, Dr. Heba Mahmoud Ali, Dr. Hend Mahmoud, Dr. Servinaz Sayed Mohammad, Dr. Rana Abo-Zaid, Dr. Veronia Rafat, Dr. Shimaa Ahmad Mekki, Dr. Amira Sobhy"
$ws.Range("G41").Value = "Dr. Alshimaa Atef, Dr. Hend Mahmoud, Dr. Mohammad El-Tanany, Dr. Hanan Ragab, Dr. Shimaa Ahmad Mekki, Dr. Amira Sobhy"
$ws.Range("G42").Value = "Dr. Menna tuâ€™Allah Medhat, Dr. Alshimaa Atef, Dr. Eman Tantawi, Dr. Servinaz Sayed Mohammad, Dr. Shimaa Ahmad Mekki"
$ws.Range("G43").Value = "Dr. Amera Ahmad Saad, Dr. Lamiaa Ossama, Dr. Fatma Elhady, Dr. Nada Mohammad, Dr. Kerelos Zareef, Dr. Menna tu'Alllah Mohammad, Dr. Abeer Ragab"
$ws.Range("G47").Value = "Dr. Maryam Ahmad, Dr. Nourhan Osama, Dr. Merna Said, Dr. Arwa Al-Sayed, Dr. Esraa Mostafa, Dr. Amira Ibrahim"
$ws.Range("G48").Value = "Dr. Maryam Ahmad, Dr. Fatma Shoukry, Dr. Amany Raafat, Dr. Merna Said, Dr. Sarah Abdelmohsen, Dr. Yasmeena Fattoh, Dr. Eman M. Abo-Sakaya"
$ws.Range("G49").Value = "Dr. Mohammad Safwat, Dr. Mariam Toma Gerges"
$ws.Range("G50").Value = "Dr. Mona Ibrahim Hussein, Dr. Heba Al-Sayed Mohammad, Dr. Dalia Tarek Elwan"
$ws.Range("G51").Value = "Menna tuâ€™Allah Gamil, Dr. Nouran Mahmoud"
$ws.Range("G56").Value = "Dr. Eman Samir Gabry, Dr. Wafaa Ebida, Dr. Ola Abd Al-Fattah, Dr. Abdullah El-Agrody"
$ws.Range("G57").Value = "Dr. Remon, Dr. Marina Atef, Dr. Nardine, Dr. Yasmin, Dr. Monica, Dr. Neveen Nashaat, Dr. Marina Sorial"
$ws.Range("G58").Value = "Dr. Remon, Dr. Marina Atef, Dr. Nardine, Dr. Yasmin, Dr. Monica, Dr. Neveen Nashaat, Dr. Marina Sorial"
$ws.Range("G59").Value = "Dr. Asmaa Reda, Dr. Nesma, Dr. Nourhan Mahmoud, Dr. Heba Mahmoud Ali, Dr. Servinaz Sayed Mohammad, Dr. Mohammad El-Tanany, Dr. Amira Sobhy"
$ws.Range("G60").Value = "Dr. Alshimaa Atef, Dr. Hend Mahmoud, Dr. Mohammad El-Tanany, Dr. Hanan Ragab, Dr. Shimaa Ahmad Mekki, Dr. Amira Sobhy"
$ws.Range("G61").Value = "Dr. Asmaa Reda, Dr. Nahla Nagiub, Dr. Majorelle Magdy, Dr. Shimaa Ahmad Mekki, Dr. Amira Sobhy"
$ws.Range("G63").Value = "Dr. Aya Saeed, Dr. Safa Hany, Dr. Amal Awwad"
$ws.Range("G66").Value = "Dr. Marina Youhanna, Dr. Dina Adel, Dr. Amira Ibrahim, Dr. Madeha Saeed, Dr. Yasmeena Fattoh, Dr. Eman M. Abo-Sakaya"
$ws.Range("G67").Value = "Dr. Esraa Mostafa, Dr. Amira Ibrahim, Dr. Yasmeena Fattoh"
$ws.Range("G75").Value = "Dr. Eman Samir Gabry, Dr. Wafaa Ebida, Dr. Ola Abd Al-Fattah, Dr. Abdullah El-Agrody"
$ws.Range("G76").Value = "Dr. Eman Samir Gabry, Dr. Yasmin, Dr. Wafaa Ebida, Dr. Neveen Nashaat, Dr. Marina Sorial"
$ws.Range("G77").Value = "Dr. Remon, Dr. Marina Atef, Dr. Nardine, Dr. Yasmin, Dr. Monica, Dr. Neveen Nashaat, Dr. Marina Sorial"
$ws.Range("G78").Value = "Dr. Asmaa Reda, Dr. Nesma, Dr. Nourhan Mahmoud, Dr. Heba Mahmoud Ali, Dr. Servinaz Sayed Mohammad, Dr. Mohammad El-Tanany, Dr. Amira Sobhy"
$ws.Range("G79").Value = "Dr. Alshimaa Atef, Dr. Hend Mahmoud, Dr. Mohammad El-Tanany, Dr. Hanan Ragab, Dr. Shimaa Ahmad Mekki, Dr. Amira Sobhy"
$ws.Range("G80").Value = "Dr. Asmaa Reda, Dr. Nahla Nagiub, Dr. Majorelle Magdy, Dr. Shimaa Ahmad Mekki, Dr. Amira Sobhy"
$ws.Range("G81").Value = "Dr. Amera Ahmad Saad, Dr. Lamiaa Ossama, Dr. Fatma Elhady, Dr. Nada Mohammad, Dr. Abeer Ragab"
$ws.Range("G83").Value = "Dr. Aya Saeed, Dr. Safa Hany, Dr. Amal Awwad"
$ws.Range("G85").Value = "Dr. Marina Youhanna, Dr. Dina Adel, Dr. Amira Ibrahim, Dr. Madeha Saeed, Dr. Yasmeena Fattoh, Dr. Eman M. Abo-Sakaya"
$ws.Range("G86").Value = "Dr. Maryam Ahmad, Dr. Fatma Shoukry, Dr. Amany Raafat, Dr. Merna Said, Dr. Sarah Abdelmohsen, Dr. Yasmeena Fattoh, Dr. Eman M. Abo-Sakaya"
$ws.Range("G88").Value = "Dr. Mona Ibrahim Hussein, Dr. Heba Al-Sayed Mohammad, Dr. Dalia Tarek Elwan"
$ws.Range("G94").Value = "Dr. Eman Samir Gabry, Dr. Wafaa Ebida, Dr. Ola Abd Al-Fattah, Dr. Abdullah El-Agrody"
$ws.Range("G95").Value = "Dr. Eman Samir Gabry, Dr. Yasmin, Dr. Wafaa Ebida, Dr. Neveen Nashaat, Dr. Marina Sorial"
$ws.Range("G96").Value = "Dr. Remon, Dr. Marina Atef, Dr. Nardine, Dr. Yasmin, Dr. Monica, Dr. Neveen Nashaat, Dr. Marina Sorial"
$ws.Range("G97").Value = "Dr. Asmaa Reda, Dr. Nesma, Dr. Nourhan Mahmoud, Dr. Heba Mahmoud Ali, Dr. Servinaz Sayed Mohammad, Dr. Mohammad El-Tanany, Dr. Amira Sobhy"
$ws.Range("G98").Value = "Dr. Alshimaa Atef, Dr. Hend Mahmoud, Dr. Mohammad El-Tanany, Dr. Hanan Ragab, Dr. Shimaa Ahmad Mekki, Dr. Amira Sobhy"
$ws.Range("G99").Value = "Dr. Menna tuâ€™Allah Medhat, Dr. Alshimaa Atef, Dr. Eman Tantawi, Dr. Servinaz Sayed Mohammad, Dr. Shimaa Ahmad Mekki"
$ws.Range("G100").Value = "Dr. Amera Ahmad Saad, Dr. Lamiaa Ossama, Dr. Fatma Elhady, Dr. Nada Mohammad, Dr. Kerelos Zareef, Dr. Menna tu'Alllah Mohammad, Dr. Abeer Ragab"
$ws.Range("G101").Value = "Dr. Aya Saeed, Dr. Safa Hany, Dr. Amal Awwad"
$ws.Range("G104").Value = "Dr. Maryam Ahmad, Dr. Nourhan Osama, Dr. Merna Said, Dr. Arwa Al-Sayed, Dr. Esraa Mostafa, Dr. Amira Ibrahim"
$ws.Range("G105").Value = "Dr. Basma Hamed, Dr. Amira Ibrahim"
$ws.Range("G113").Value = "Dr. Eman Samir Gabry, Dr. Wafaa Ebida, Dr. Ola Abd Al-Fattah, Dr. Abdullah El-Agrody"
$ws.Range("G115").Value = "Dr. Remon, Dr. Marina Atef, Dr. Nardine, Dr. Yasmin, Dr. Monica, Dr. Neveen Nashaat, Dr. Marina Sorial"
